# "Current Account - Service Domain Specification.xlsx" update
#
# The six separate "Key Features" bullet rows (previously rows 9-13, plus
# the first bullet already living in B8) are consolidated into a single
# B8 cell holding one combined sentence. The now-empty bullet rows 9-13
# are removed, which shifts the "Functional Pattern" / "Asset Type" /
# "Generic Artifact" / "Control Record" / "Registration Status" block
# (previously rows 14-18) up to rows 9-13. The named range and sheet
# dimension shrink accordingly from row 18 to row 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge all the "Key Features" bullets into a single cell (B8), replacing
# the first bullet ("Set-up current account facilities") that used to sit
# there alone.
$ws.Range("B8").Value = "Set-up current account facilities, Issue cards, check-books for an account, Process deposits and withdrawals, Handle standing orders and direct debits.Process fees and apply interest charges. Provide balance/statements and reports."

# The wrapped text now needs more vertical room.
$ws.Rows.Item(8).RowHeight = 43.2

# The other five bullet rows (formerly rows 9-13, column A only) are no
# longer needed now that their text lives inside B8 - remove them. This
# pulls the rows below (old 14-18) up to become the new rows 9-13.
$ws.Range("A9:B13").EntireRow.Delete()

# Update the workbook-level named range so it still spans the full table,
# which now ends at row 13 instead of row 18.
$name = $wb.Names.Item("CurrentAccountSpecification")
$name.RefersTo = "=Sheet1!`$A`$1:`$B`$13"

# Reflect the saved cursor position from the authored workbook.
$ws.Range("B26").Select()
